$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header cells for 2023 and 2024 (same visual style family as D1's header)
$ws.Range("E1").Value = "Cantidad de víctimas (enero - diciembre 2023)"
$ws.Range("F1").Value = "Cantidad de víctimas (enero - diciembre 2024)"

# New data values for rows 2-5
$ws.Range("E2").Value = 984
$ws.Range("F2").Value = 1369

$ws.Range("E3").Value = 33
$ws.Range("F3").Value = 94

$ws.Range("E4").Value = 0
$ws.Range("F4").Value = 0

$ws.Range("E5").Value = 21
$ws.Range("F5").Value = 32

# Widen D:F to fit the new, longer header text (matches the ~39.25-char
# "best fit" width Excel settles on for these headers)
$ws.Range("D1:F1").EntireColumn.ColumnWidth = 38.41666667

# Update selection to reflect the new active cell
$ws.Range("F12").Select()
